# The document's headers/footers each contain one inline picture (the
# BTec logo in the headers, the Pearson logo in the footers). This edit
# swaps which generic "imageN.ext" display name is reported for each of
# those inline pictures:
#   - header pictures (BTec_Logo-Orange):                image2.jpg -> image1.jpg
#   - footer pictures (...\PearsonLogo.png, alt text):    image1.png -> image2.png
#
# The pictures themselves (their embedded binary / relationship target)
# are untouched - only the shape's reported Name changes.
#
# Note: InlineShape.Name doesn't round-trip as a readable value in this
# host, so rather than toggling whatever the current name happens to be,
# the picture is identified by its (reliably readable) AlternativeText /
# its header-vs-footer story, and the name is set to its known target
# value directly.

$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hfIndex in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {

        foreach ($hfColl in @($sec.Headers, $sec.Footers)) {
            $hf = $hfColl.Item($hfIndex)
            if ($hf.Exists) {
                $shapes = $hf.Range.InlineShapes
                for ($i = 1; $i -le $shapes.Count; $i++) {
                    $shp = $shapes.Item($i)
                    if ($hf.IsHeader) {
                        $shp.Name = "image1.jpg"
                    } else {
                        $shp.Name = "image2.png"
                    }
                }
            }
        }
    }
}
